$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.474661827087402
$ws.Range("B1").Value = 3.787441253662109
$ws.Range("C1").Value = 3.73789381980896
$ws.Range("D1").Value = 1.572843074798584
$ws.Range("E1").Value = 0.9992862343788147
